# Scheduled runner data refresh: update market price / profit columns (H-N)
# across all job sheets per the latest Universalis snapshot.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 121.46667
$ws.Range("I6").Value = 108.9
$ws.Range("J6").Value = 146.6
$ws.Range("K6").Value = 326.7
$ws.Range("L6").Value = 439.8
$ws.Range("M6").Value = -214.7
$ws.Range("N6").Value = -663.8
$ws.Range("H18").Value = 2490.5
$ws.Range("I18").Value = 2490.5
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 2490.5
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -2206.5
$ws.Range("N18").ClearContents()
$ws.Range("H48").Value = 808.5
$ws.Range("I48").Value = 808.5
$ws.Range("K48").Value = 2425.5
$ws.Range("M48").Value = -2133.5
$ws.Range("H56").Value = 808.5
$ws.Range("I56").Value = 808.5
$ws.Range("K56").Value = 2425.5
$ws.Range("M56").Value = -1891.5
$ws.Range("H86").Value = 2544.261
$ws.Range("I86").Value = 2591.7273
$ws.Range("J86").Value = 1500
$ws.Range("K86").Value = 2591.7273
$ws.Range("L86").Value = 1500
$ws.Range("M86").Value = -1468.7273
$ws.Range("N86").Value = -3746
$ws.Range("H89").Value = 2544.261
$ws.Range("I89").Value = 2591.7273
$ws.Range("J89").Value = 1500
$ws.Range("K89").Value = 12958.6365
$ws.Range("L89").Value = 7500
$ws.Range("M89").Value = -7342.636500000001
$ws.Range("N89").Value = -18732
$ws.Range("H99").Value = 4896
$ws.Range("J99").Value = 6455
$ws.Range("L99").Value = 19365
$ws.Range("N99").Value = -22361
$ws.Range("H101").Value = 895.63635
$ws.Range("J101").Value = 1086.6428
$ws.Range("L101").Value = 3259.9284
$ws.Range("N101").Value = -6503.928400000001
$ws.Range("H129").Value = 9068.777
$ws.Range("I129").Value = 1012.5
$ws.Range("K129").Value = 3037.5
$ws.Range("M129").Value = 1962.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 806.9048
$ws.Range("J2").Value = 1499.8
$ws.Range("L2").Value = 1499.8
$ws.Range("N2").Value = -1725.8
$ws.Range("H21").Value = 45398
$ws.Range("I21").Value = 1990
$ws.Range("J21").Value = 56250
$ws.Range("K21").Value = 1990
$ws.Range("L21").Value = 56250
$ws.Range("M21").Value = -1616
$ws.Range("N21").Value = -56998
$ws.Range("H61").Value = 9553566
$ws.Range("I61").Value = 10534708
$ws.Range("J61").Value = 3339666.2
$ws.Range("K61").Value = 10534708
$ws.Range("L61").Value = 3339666.2
$ws.Range("M61").Value = -10534496
$ws.Range("N61").Value = -3340090.2
$ws.Range("H116").Value = 806.9048
$ws.Range("J116").Value = 1499.8
$ws.Range("L116").Value = 1499.8
$ws.Range("N116").Value = -6087.8
$ws.Range("H132").Value = 4352376
$ws.Range("I132").Value = 4928.8423
$ws.Range("K132").Value = 14786.5269
$ws.Range("M132").Value = -12256.5269
$ws.Range("H136").Value = 9553566
$ws.Range("I136").Value = 10534708
$ws.Range("J136").Value = 3339666.2
$ws.Range("K136").Value = 31604124
$ws.Range("L136").Value = 10018998.6
$ws.Range("M136").Value = -31601574
$ws.Range("N136").Value = -10024098.6

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 806.9048
$ws.Range("J3").Value = 1499.8
$ws.Range("L3").Value = 1499.8
$ws.Range("N3").Value = -1727.8
$ws.Range("H22").Value = 1047.875
$ws.Range("I22").Value = 1253.5
$ws.Range("J22").Value = 431
$ws.Range("K22").Value = 1253.5
$ws.Range("L22").Value = 431
$ws.Range("M22").Value = -1080.5
$ws.Range("N22").Value = -777
$ws.Range("H57").Value = 89153.78
$ws.Range("J57").Value = 89153.78
$ws.Range("L57").Value = 89153.78
$ws.Range("N57").Value = -90593.78
$ws.Range("H134").Value = 100000000
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 100000000
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 300000000
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -300005070
$ws.Range("H136").Value = 89153.78
$ws.Range("J136").Value = 89153.78
$ws.Range("L136").Value = 89153.78
$ws.Range("N136").Value = -99353.78

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H54").Value = 49999
$ws.Range("J54").Value = 49999
$ws.Range("L54").Value = 49999
$ws.Range("N54").Value = -51315
$ws.Range("H94").Value = 1999.4286
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 1999.4286
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 1999.4286
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -2901.4286
$ws.Range("H141").Value = 415363.47
$ws.Range("J141").Value = 474333.12
$ws.Range("L141").Value = 474333.12
$ws.Range("N141").Value = -484693.12

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9383.691999999999
$ws.Range("I3").Value = 7241.364
$ws.Range("J3").Value = 21166.5
$ws.Range("K3").Value = 21724.092
$ws.Range("L3").Value = 63499.5
$ws.Range("M3").Value = -21612.092
$ws.Range("N3").Value = -63723.5
$ws.Range("H62").Value = 12172.625
$ws.Range("I62").Value = 8690
$ws.Range("K62").Value = 26070
$ws.Range("M62").Value = -25384
$ws.Range("H65").Value = 12172.625
$ws.Range("I65").Value = 8690
$ws.Range("K65").Value = 78210
$ws.Range("M65").Value = -74778
$ws.Range("H134").Value = 8174.15
$ws.Range("I134").Value = 2509.375
$ws.Range("K134").Value = 7528.125
$ws.Range("M134").Value = -2458.125
$ws.Range("H136").Value = 1843
$ws.Range("I136").Value = 1211.6
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 3634.8
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = 1465.2
$ws.Range("N136").Value = -25200
$ws.Range("H140").Value = 2897.487
$ws.Range("I140").Value = 1184.303
$ws.Range("J140").Value = 12320
$ws.Range("K140").Value = 3552.909000000001
$ws.Range("L140").Value = 36960
$ws.Range("M140").Value = 1627.090999999999
$ws.Range("N140").Value = -47320
$ws.Range("H141").Value = 4121.8335
$ws.Range("I141").Value = 4121.8335
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 12365.5005
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -7185.500499999998
$ws.Range("N141").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 10578.4
$ws.Range("I22").Value = 10578.4
$ws.Range("K22").Value = 10578.4
$ws.Range("M22").Value = -10049.4
$ws.Range("H54").Value = 24666.334
$ws.Range("J54").Value = 24666.334
$ws.Range("L54").Value = 24666.334
$ws.Range("N54").Value = -25446.334
$ws.Range("H134").Value = 63332.832
$ws.Range("J134").Value = 63332.832
$ws.Range("L134").Value = 189998.496
$ws.Range("N134").Value = -195068.496

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 29062.143
$ws.Range("I3").Value = 14816.667
$ws.Range("K3").Value = 14816.667
$ws.Range("M3").Value = -14704.667
$ws.Range("H15").Value = 29062.143
$ws.Range("I15").Value = 14816.667
$ws.Range("K15").Value = 14816.667
$ws.Range("M15").Value = -14646.667
$ws.Range("H17").Value = 13436
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 13436
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 13436
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -13776
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("H31").Value = 4578.4287
$ws.Range("I31").Value = 4008.3333
$ws.Range("J31").Value = 7999
$ws.Range("K31").Value = 4008.3333
$ws.Range("L31").Value = 7999
$ws.Range("M31").Value = -3760.3333
$ws.Range("N31").Value = -8495
$ws.Range("H137").Value = 116280.86
$ws.Range("J137").Value = 116280.86
$ws.Range("L137").Value = 116280.86
$ws.Range("N137").Value = -126480.86

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 202020.6
$ws.Range("I132").Value = 1805.5814
$ws.Range("K132").Value = 5416.7442
$ws.Range("M132").Value = -2886.7442

